$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting id/name/cells/kwargs to the right
$ws.Range("B1").EntireColumn.Insert()

# Set the new header value, copying the formatting of the header row
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B1").Value = "env"
